$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1 (time changed from 04:02 to 05:19)
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 05:19"

# Row 55 (Honduras) - refreshed case counts
$ws.Range("B55").Value = 20262
$ws.Range("C55").Value = 704
$ws.Range("D55").Value = 2123
$ws.Range("E55").Value = 17597
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 45
$ws.Range("H55").Value = 542

# Haiti moves above Tayikistan in the country ranking, with refreshed data.
# Row 83 becomes Haiti (new figures), row 84 becomes Tayikistan (its former,
# unchanged figures that used to sit in row 83).
$ws.Range("A83").Value = "Haiti"
$ws.Range("B83").Value = 6021
$ws.Range("C83").Value = 46
$ws.Range("D83").Value = 931
$ws.Range("E83").Value = 4983
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 107

$ws.Range("A84").Value = "Tayikistan"
$ws.Range("B84").Value = 6005
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 4627
$ws.Range("E84").Value = 1326
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 52

# Laos moves above Santa Lucia in the ranking (figures identical, only the
# country names/order swap).
$ws.Range("A203").Value = "Laos"
$ws.Range("A204").Value = "Santa Lucia"
